# Faturamento diario - atualizei dados bibi e add
# Insert 2 new daily rows into the June/2025 block (after day 27, row 28)
# and correct the "day 27" total for June.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing June/27 total_venda value.
$ws.Range("B28").Value = 24519.52

# Make room for two new days (28 and 29 of June/2025); everything below
# shifts down by two rows automatically.
$ws.Rows("29:30").Insert()

# Day 28 of June/2025.
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = 10903.9
$ws.Range("C29").Value = 6
$ws.Range("D29").Value = 2025
$ws.Range("E29").Value = "06/2025"

# Day 29 of June/2025.
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = 6129.8
$ws.Range("C30").Value = 6
$ws.Range("D30").Value = 2025
$ws.Range("E30").Value = "06/2025"
